$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.622153
$ws.Range("H2").Value = 1.866459
$ws.Range("I2").Value = 0.5959696685805808
$ws.Range("J2").Value = 0.5959696685805808
$ws.Range("M2").Value = 2.906846333333333
$ws.Range("N2").Value = 8.720538999999999
$ws.Range("O2").Value = 0.005520525738044089
$ws.Range("P2").Value = 0.005624540846623205
$ws.Range("Q2").Value = 1.808503166822333
$ws.Range("R2").Value = 16.276528501401
$ws.Range("S2").Value = 0.003290065894492701
$ws.Range("T2").Value = 0.003352055744279971
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.622153
$ws.Range("H3").Value = 1.866459
$ws.Range("I3").Value = 0.5959696685805808
$ws.Range("J3").Value = 0.5959696685805808
$ws.Range("O3").Value = 0.3528665483720876
$ws.Range("P3").Value = 0.3595150912979765
$ws.Range("Q3").Value = 115.5977347952163
$ws.Range("R3").Value = 1040.379613156947
$ws.Range("S3").Value = 0.2102977598864866
$ws.Range("T3").Value = 0.2142600898105723
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.622153
$ws.Range("H4").Value = 1.866459
$ws.Range("I4").Value = 0.5959696685805808
$ws.Range("J4").Value = 0.5959696685805808
$ws.Range("M4").Value = 137.0717086666666
$ws.Range("N4").Value = 411.2151259999999
$ws.Range("O4").Value = 0.2603191943704447
$ws.Range("P4").Value = 0.2652240042658267
$ws.Range("Q4").Value = 85.27957476209264
$ws.Range("R4").Value = 767.5161728588338
$ws.Range("S4").Value = 0.1551423439941177
$ws.Range("T4").Value = 0.1580654619219193
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.622153
$ws.Range("H5").Value = 1.866459
$ws.Range("I5").Value = 0.5959696685805808
$ws.Range("J5").Value = 0.5959696685805808
$ws.Range("M5").Value = 29.2127365
$ws.Range("N5").Value = 58.425473
$ws.Range("O5").Value = 0.05547925319534149
$ws.Range("P5").Value = 0.03768304451958546
$ws.Range("Q5").Value = 18.1747916516845
$ws.Range("R5").Value = 109.048749910107
$ws.Range("S5").Value = 0.03306395213992579
$ws.Range("T5").Value = 0.02245795155344461
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.622153
$ws.Range("H6").Value = 1.866459
$ws.Range("I6").Value = 0.5959696685805808
$ws.Range("J6").Value = 0.5959696685805808
$ws.Range("M6").Value = 171.5584106666666
$ws.Range("N6").Value = 514.6752319999999
$ws.Range("O6").Value = 0.3258144783240821
$ws.Range("P6").Value = 0.331953319069988
$ws.Range("Q6").Value = 106.7355798714986
$ws.Range("R6").Value = 960.6202188434878
$ws.Range("S6").Value = 0.194175546665558
$ws.Range("T6").Value = 0.1978341095503645
$ws.Range("G7").Value = 0.421781
$ws.Range("H7").Value = 1.265343
$ws.Range("I7").Value = 0.4040303314194192
$ws.Range("J7").Value = 0.4040303314194192
$ws.Range("M7").Value = 2.906846333333333
$ws.Range("N7").Value = 8.720538999999999
$ws.Range("O7").Value = 0.005520525738044089
$ws.Range("P7").Value = 0.005624540846623205
$ws.Range("Q7").Value = 1.226052553319666
$ws.Range("R7").Value = 11.034472979877
$ws.Range("S7").Value = 0.002230459843551387
$ws.Range("T7").Value = 0.002272485102343234
$ws.Range("G8").Value = 0.421781
$ws.Range("H8").Value = 1.265343
$ws.Range("I8").Value = 0.4040303314194192
$ws.Range("J8").Value = 0.4040303314194192
$ws.Range("O8").Value = 0.3528665483720876
$ws.Range("P8").Value = 0.3595150912979765
$ws.Range("Q8").Value = 78.36806730765765
$ws.Range("R8").Value = 705.312605768919
$ws.Range("S8").Value = 0.1425687884856011
$ws.Range("T8").Value = 0.1452550014874042
$ws.Range("G9").Value = 0.421781
$ws.Range("H9").Value = 1.265343
$ws.Range("I9").Value = 0.4040303314194192
$ws.Range("J9").Value = 0.4040303314194192
$ws.Range("M9").Value = 137.0717086666666
$ws.Range("N9").Value = 411.2151259999999
$ws.Range("O9").Value = 0.2603191943704447
$ws.Range("P9").Value = 0.2652240042658267
$ws.Range("Q9").Value = 57.81424235313531
$ws.Range("R9").Value = 520.3281811782178
$ws.Range("S9").Value = 0.105176850376327
$ws.Range("T9").Value = 0.1071585423439074
$ws.Range("G10").Value = 0.421781
$ws.Range("H10").Value = 1.265343
$ws.Range("I10").Value = 0.4040303314194192
$ws.Range("J10").Value = 0.4040303314194192
$ws.Range("M10").Value = 29.2127365
$ws.Range("N10").Value = 58.425473
$ws.Range("O10").Value = 0.05547925319534149
$ws.Range("P10").Value = 0.03768304451958546
$ws.Range("Q10").Value = 12.3213772137065
$ws.Range("R10").Value = 73.92826328223899
$ws.Range("S10").Value = 0.02241530105541569
$ws.Range("T10").Value = 0.01522509296614084
$ws.Range("G11").Value = 0.421781
$ws.Range("H11").Value = 1.265343
$ws.Range("I11").Value = 0.4040303314194192
$ws.Range("J11").Value = 0.4040303314194192
$ws.Range("M11").Value = 171.5584106666666
$ws.Range("N11").Value = 514.6752319999999
$ws.Range("O11").Value = 0.3258144783240821
$ws.Range("P11").Value = 0.331953319069988
$ws.Range("Q11").Value = 72.36007800939731
$ws.Range("R11").Value = 651.2407020845759
$ws.Range("S11").Value = 0.1316389316585241
$ws.Range("T11").Value = 0.1341192095196235
